# ---------------------------------------------------------------------------
# 北京-漫展信息.xlsx  -  scraper refresh ("Update gh-pages to output
# generated at 456a3b4")
#
# The scrape bumps column F ("想去人数" / want-to-go count) for a batch of
# events. The same event can appear on more than one tab (展览/演出/本地生活
# each roll up into 全部类型), so the same counter is bumped on every tab
# that lists it. Cells are addressed by their A1 ref on each sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3649
$ws.Range("F3").Value = 557
$ws.Range("F4").Value = 202
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 47
$ws.Range("F8").Value = 11
$ws.Range("F9").Value = 369
$ws.Range("F10").Value = 4434
$ws.Range("F11").Value = 4434
$ws.Range("F14").Value = 1060
$ws.Range("F15").Value = 571
$ws.Range("F16").Value = 3871
$ws.Range("F17").Value = 141
$ws.Range("F18").Value = 133
$ws.Range("F20").Value = 146
$ws.Range("F21").Value = 3262
$ws.Range("F24").Value = 11
$ws.Range("F25").Value = 2709
$ws.Range("F26").Value = 97
$ws.Range("F27").Value = 105
$ws.Range("F29").Value = 120
$ws.Range("F31").Value = 152
$ws.Range("F32").Value = 62
$ws.Range("F33").Value = 34
$ws.Range("F34").Value = 17
$ws.Range("F35").Value = 47
$ws.Range("F36").Value = 121
$ws.Range("F37").Value = 4978
$ws.Range("F38").Value = 652
$ws.Range("F39").Value = 367
$ws.Range("F40").Value = 72
$ws.Range("F41").Value = 948
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 956
$ws.Range("F44").Value = 354
$ws.Range("F46").Value = 1838
$ws.Range("F47").Value = 282
$ws.Range("F49").Value = 660
$ws.Range("F50").Value = 788

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 7
$ws.Range("F6").Value = 74
$ws.Range("F15").Value = 117
$ws.Range("F22").Value = 692

# --- 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 187
$ws.Range("F3").Value = 12

# --- 全部类型 (All types - aggregates the other three tabs) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 187
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 3649
$ws.Range("F5").Value = 3649
$ws.Range("F7").Value = 557
$ws.Range("F8").Value = 202
$ws.Range("F9").Value = 10
$ws.Range("F10").Value = 74
$ws.Range("F11").Value = 47
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = 369
$ws.Range("F14").Value = 4435
$ws.Range("F15").Value = 4435
$ws.Range("F20").Value = 1060
$ws.Range("F21").Value = 571
$ws.Range("F22").Value = 3871
$ws.Range("F23").Value = 141
$ws.Range("F24").Value = 133
$ws.Range("F25").Value = 3262
$ws.Range("F26").Value = 2709
$ws.Range("F27").Value = 97
$ws.Range("F28").Value = 105
$ws.Range("F29").Value = 120
$ws.Range("F31").Value = 152
$ws.Range("F32").Value = 62
$ws.Range("F33").Value = 17
$ws.Range("F35").Value = 47
$ws.Range("F36").Value = 121
$ws.Range("F37").Value = 117
$ws.Range("F39").Value = 4978
$ws.Range("F41").Value = 652
$ws.Range("F42").Value = 367
$ws.Range("F44").Value = 72
$ws.Range("F45").Value = 956
$ws.Range("F46").Value = 354
$ws.Range("F48").Value = 1838
$ws.Range("F49").Value = 282
$ws.Range("F51").Value = 788

